$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'25.834.99"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -2.57%  '

$ws.Cells.Item(3, 4).Value = "'1.751.37"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -5.01%  '

$ws.Cells.Item(4, 4).Value = "'0.9950"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.59%  '

$ws.Cells.Item(5, 4).Value = "'238.92"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -8.98%  '

$ws.Cells.Item(6, 4).Value = "'0.9963"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.46%  '

$ws.Cells.Item(7, 4).Value = "'0.4984"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -6.49%  '

$ws.Cells.Item(8, 4).Value = "'41.67"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -7.39%  '

$ws.Cells.Item(9, 4).Value = "'0.2661"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -13.44%  '

$ws.Cells.Item(10, 4).Value = "'0.06173"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -10.52%  '

$ws.Cells.Item(11, 4).Value = "'1.744.05"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -5.38%  '

$ws.Cells.Item(12, 4).Value = "'0.06932"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -11.28%  '

$ws.Cells.Item(13, 4).Value = "'15.52"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -15.73%  '

$ws.Cells.Item(14, 4).Value = "'4.489"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -10.48%  '

$ws.Cells.Item(15, 2).Value = 'Polygon'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(15, 4).Value = "'0.5929"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -21.36%  '

$ws.Cells.Item(16, 2).Value = 'Litecoin'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(16, 4).Value = "'77.13"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -14.12%  '

$ws.Cells.Item(17, 4).Value = "'0.9966"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.37%  '

$ws.Cells.Item(18, 4).Value = "'0.9971"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -0.39%  '

$ws.Cells.Item(19, 4).Value = "'25.839.41"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -2.60%  '

$ws.Cells.Item(20, 4).Value = "'11.70"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -16.48%  '

$ws.Cells.Item(21, 4).Value = "'0.000006734"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -15.30%  '

$ws.Cells.Item(22, 4).Value = "'1.963.49"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -5.57%  '

$ws.Cells.Item(23, 4).Value = "'4.069"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -11.97%  '

$ws.Cells.Item(24, 4).Value = "'8.038"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -13.83%  '

$ws.Cells.Item(25, 4).Value = "'5.105"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -14.75%  '

$ws.Cells.Item(26, 4).Value = "'137.56"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -3.69%  '

$ws.Cells.Item(27, 4).Value = "'1.525"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -9.68%  '

$ws.Cells.Item(28, 4).Value = "'1.854"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -15.49%  '

$ws.Cells.Item(29, 4).Value = "'14.90"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -12.43%  '

$ws.Cells.Item(30, 4).Value = "'102.49"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -7.90%  '

$ws.Cells.Item(31, 4).Value = "'3.766"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -11.89%  '

$ws.Cells.Item(32, 4).Value = "'0.08053"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -8.60%  '

$ws.Cells.Item(33, 4).Value = "'3.478"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -14.86%  '

$ws.Cells.Item(34, 4).Value = "'0.04483"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -6.90%  '

$ws.Cells.Item(35, 4).Value = "'0.9961"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -0.41%  '

$ws.Cells.Item(36, 4).Value = "'2.615"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -10.74%  '

$ws.Cells.Item(37, 4).Value = "'0.9893"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -12.86%  '

$ws.Cells.Item(38, 4).Value = "'0.6055"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -17.19%  '

$ws.Cells.Item(39, 4).Value = "'2.663"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -14.30%  '

$ws.Cells.Item(40, 2).Value = 'RenderToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(40, 4).Value = "'1.952"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -15.32%  '

$ws.Cells.Item(41, 2).Value = 'Quant'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(41, 4).Value = "'105.21"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -2.90%  '

$ws.Cells.Item(42, 2).Value = 'VeChain'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(42, 4).Value = "'0.01522"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -11.37%  '

$ws.Cells.Item(43, 2).Value = 'PaxDollar'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(43, 4).Value = "'0.9970"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -0.40%  '

$ws.Cells.Item(44, 2).Value = 'TheSandbox'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(44, 4).Value = "'0.3825"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -20.27%  '

$ws.Cells.Item(45, 2).Value = 'FraxShare'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(45, 4).Value = "'5.161"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -12.01%  '

$ws.Cells.Item(46, 4).Value = "'0.7306"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -19.47%  '

$ws.Cells.Item(47, 4).Value = "'0.1117"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -10.18%  '

$ws.Cells.Item(48, 4).Value = "'0.05225"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -9.85%  '

$ws.Cells.Item(49, 2).Value = 'Aptos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(49, 4).Value = "'5.975"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -20.51%  '

$ws.Cells.Item(50, 2).Value = 'Elrond'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Cells.Item(50, 4).Value = "'30.25"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -13.32%  '

$ws.Cells.Item(51, 4).Value = "'52.54"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -12.82%  '
